# Update cryptocurrency price (D) and 1h volume change (E) columns
# with freshly scraped values from coinranking.com
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new text would otherwise be auto-parsed as a
# number by Excel (single "." and all-digits) are pre-formatted as Text
# so they are stored/display exactly like the original inline strings.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Column D (Price) updates
$ws.Range("D2").Value = "63.519.24"
$ws.Range("D3").Value = "2.604.05"
$ws.Range("D5").Value = "588.06"
$ws.Range("D6").Value = "149.34"
$ws.Range("D13").Value = "27.55"
$ws.Range("D14").Value = "3.071.20"
$ws.Range("D15").Value = "63.330.10"
$ws.Range("D17").Value = "2.590.53"
$ws.Range("D20").Value = "344.20"
$ws.Range("D21").Value = "6.84"
$ws.Range("D23").Value = "66.44"
$ws.Range("D24").Value = "1.73"
$ws.Range("D27").Value = "564.41"
$ws.Range("D30").Value = "1.00"
$ws.Range("D32").Value = "0.0₃0845"
$ws.Range("D35").Value = "165.73"
$ws.Range("D37").Value = "0.999"
$ws.Range("D38").Value = "19.42"
$ws.Range("D41").Value = "165.71"
$ws.Range("D43").Value = "22.77"
$ws.Range("D49").Value = "19.09"

# Column E (Volume 1h %) updates
$ws.Range("E2").Value = "  -1.42%  "
$ws.Range("E3").Value = "  -1.55%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -3.11%  "
$ws.Range("E6").Value = "  -1.69%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -1.75%  "
$ws.Range("E9").Value = "  -1.36%  "
$ws.Range("E10").Value = "  +0.91%  "
$ws.Range("E11").Value = "  -1.17%  "
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("E13").Value = "  -1.31%  "
$ws.Range("E14").Value = "  -1.55%  "
$ws.Range("E15").Value = "  -1.43%  "
$ws.Range("E16").Value = "  +3.53%  "
$ws.Range("E17").Value = "  -2.29%  "
$ws.Range("E18").Value = "  -1.59%  "
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("E20").Value = "  -2.73%  "
$ws.Range("E21").Value = "  -2.11%  "
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("E23").Value = "  -0.68%  "
$ws.Range("E24").Value = "  -1.98%  "
$ws.Range("E25").Value = "  -2.17%  "
$ws.Range("E26").Value = "  -4.43%  "
$ws.Range("E27").Value = "  +2.34%  "
$ws.Range("E28").Value = "  -0.59%  "
$ws.Range("E29").Value = "  -3.31%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  -3.00%  "
$ws.Range("E32").Value = "  -3.37%  "
$ws.Range("E33").Value = "  -1.57%  "
$ws.Range("E34").Value = "  -1.13%  "
$ws.Range("E35").Value = "  -1.10%  "
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("E38").Value = "  -1.23%  "
$ws.Range("E39").Value = "  -6.19%  "
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("E41").Value = "  -1.81%  "
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("E43").Value = "  +4.38%  "
$ws.Range("E44").Value = "  -1.86%  "
$ws.Range("E45").Value = "  +1.93%  "
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("E48").Value = "  -1.14%  "
$ws.Range("E49").Value = "  -2.18%  "
$ws.Range("E50").Value = "  +13.38%  "
$ws.Range("E51").Value = "  -4.12%  "
